$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Sprint Backlog" tracking columns (G = Week/Day 3 remaining, H = Week/Day 4 remaining)
# for the rows that now have actuals recorded.

$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0

$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0

$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0

$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0

$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0

$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0

# Reflect the current selection as last used by the author
$ws.Range("G25").Select()

$wb.Save()
